$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5710.7144
$ws.Range("I86").Value = 4500
$ws.Range("J86").Value = 6618.75
$ws.Range("K86").Value = 4500
$ws.Range("L86").Value = 6618.75
$ws.Range("M86").Value = -3377
$ws.Range("N86").Value = -8864.75
$ws.Range("H89").Value = 5710.7144
$ws.Range("I89").Value = 4500
$ws.Range("J89").Value = 6618.75
$ws.Range("K89").Value = 22500
$ws.Range("L89").Value = 33093.75
$ws.Range("M89").Value = -16884
$ws.Range("N89").Value = -44325.75
$ws.Range("H137").Value = 4253.6924
$ws.Range("J137").Value = 4750
$ws.Range("L137").Value = 14250
$ws.Range("N137").Value = -19350
$ws.Range("H138").Value = 3392.4153
$ws.Range("J138").Value = 3718.1458
$ws.Range("L138").Value = 11154.4374
$ws.Range("N138").Value = -21434.4374

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5088.478
$ws.Range("I61").Value = 4998.154
$ws.Range("K61").Value = 4998.154
$ws.Range("M61").Value = -4786.154
$ws.Range("H74").Value = 20835950
$ws.Range("I74").Value = 25642228
$ws.Range("J74").Value = 8754.666999999999
$ws.Range("K74").Value = 25642228
$ws.Range("L74").Value = 8754.666999999999
$ws.Range("M74").Value = -25641354
$ws.Range("N74").Value = -10502.667
$ws.Range("H77").Value = 20835950
$ws.Range("I77").Value = 25642228
$ws.Range("J77").Value = 8754.666999999999
$ws.Range("K77").Value = 128211140
$ws.Range("L77").Value = 43773.335
$ws.Range("M77").Value = -128206772
$ws.Range("N77").Value = -52509.335
$ws.Range("H122").Value = 3520.05
$ws.Range("J122").Value = 4750
$ws.Range("L122").Value = 14250
$ws.Range("N122").Value = -19150
$ws.Range("H132").Value = 4450.9414
$ws.Range("I132").Value = 3158.6667
$ws.Range("J132").Value = 5904.75
$ws.Range("K132").Value = 9476.000100000001
$ws.Range("L132").Value = 17714.25
$ws.Range("M132").Value = -6946.000100000001
$ws.Range("N132").Value = -22774.25
$ws.Range("H136").Value = 5088.478
$ws.Range("I136").Value = 4998.154
$ws.Range("K136").Value = 14994.462
$ws.Range("M136").Value = -12444.462

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 13170.448
$ws.Range("I105").Value = 11367.523
$ws.Range("J105").Value = 17903.125
$ws.Range("K105").Value = 11367.523
$ws.Range("L105").Value = 17903.125
$ws.Range("M105").Value = -9620.522999999999
$ws.Range("N105").Value = -21397.125
$ws.Range("H134").Value = 3256.2173
$ws.Range("I134").Value = 1836.0714
$ws.Range("J134").Value = 5465.3335
$ws.Range("K134").Value = 5508.2142
$ws.Range("L134").Value = 16396.0005
$ws.Range("M134").Value = -2973.2142
$ws.Range("N134").Value = -21466.0005

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2086.111
$ws.Range("I16").Value = 1436.6666
$ws.Range("K16").Value = 1436.6666
$ws.Range("M16").Value = -1149.6666
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H58").Value = 4222.2104
$ws.Range("I58").Value = 1731.4286
$ws.Range("J58").Value = 11196.4
$ws.Range("K58").Value = 1731.4286
$ws.Range("L58").Value = 11196.4
$ws.Range("M58").Value = -1528.4286
$ws.Range("N58").Value = -11602.4
$ws.Range("H99").Value = 2281.25
$ws.Range("J99").Value = 2750
$ws.Range("L99").Value = 2750
$ws.Range("N99").Value = -5746
$ws.Range("H113").Value = 2086.111
$ws.Range("I113").Value = 1436.6666
$ws.Range("K113").Value = 1436.6666
$ws.Range("M113").Value = 733.3334
$ws.Range("H122").Value = 4109.95
$ws.Range("I122").Value = 1070.5294
$ws.Range("J122").Value = 21333.334
$ws.Range("K122").Value = 3211.5882
$ws.Range("L122").Value = 64000.00199999999
$ws.Range("M122").Value = -761.5881999999997
$ws.Range("N122").Value = -68900.00199999999
$ws.Range("H126").Value = 2281.25
$ws.Range("J126").Value = 2750
$ws.Range("L126").Value = 8250
$ws.Range("N126").Value = -13190
$ws.Range("H132").Value = 2962.6667
$ws.Range("I132").Value = 2156.375
$ws.Range("K132").Value = 6469.125
$ws.Range("M132").Value = -3939.125
$ws.Range("H134").Value = 3247.4
$ws.Range("I134").Value = 1653.7778
$ws.Range("K134").Value = 4961.3334
$ws.Range("M134").Value = -2426.3334
$ws.Range("H136").Value = 4222.2104
$ws.Range("I136").Value = 1731.4286
$ws.Range("J136").Value = 11196.4
$ws.Range("K136").Value = 5194.2858
$ws.Range("L136").Value = 33589.2
$ws.Range("M136").Value = -2644.2858
$ws.Range("N136").Value = -38689.2

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 3789.7144
$ws.Range("J126").Value = 4339.6
$ws.Range("L126").Value = 13018.8
$ws.Range("N126").Value = -22898.8
$ws.Range("H131").Value = 6946157.5
$ws.Range("I131").Value = 25000666
$ws.Range("J131").Value = 4744388
$ws.Range("K131").Value = 75001998
$ws.Range("L131").Value = 14233164
$ws.Range("M131").Value = -74996958
$ws.Range("N131").Value = -14243244

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 298441.75
$ws.Range("I80").Value = 456772.62
$ws.Range("J80").Value = 8168.5
$ws.Range("K80").Value = 456772.62
$ws.Range("L80").Value = 8168.5
$ws.Range("M80").Value = -455774.62
$ws.Range("N80").Value = -10164.5
$ws.Range("H83").Value = 298441.75
$ws.Range("I83").Value = 456772.62
$ws.Range("J83").Value = 8168.5
$ws.Range("K83").Value = 2283863.1
$ws.Range("L83").Value = 40842.5
$ws.Range("M83").Value = -2278871.1
$ws.Range("N83").Value = -50826.5
$ws.Range("H122").Value = 57503.5
$ws.Range("I122").Value = 99999
$ws.Range("J122").Value = 15008
$ws.Range("K122").Value = 299997
$ws.Range("L122").Value = 45024
$ws.Range("M122").Value = -297547
$ws.Range("N122").Value = -49924
$ws.Range("H132").Value = 3197.7942
$ws.Range("I132").Value = 2757.5356
$ws.Range("K132").Value = 8272.606800000001
$ws.Range("M132").Value = -5742.606800000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10216.137
$ws.Range("I40").Value = 9646.666999999999
$ws.Range("K40").Value = 9646.666999999999
$ws.Range("M40").Value = -9510.666999999999
$ws.Range("H55").Value = 1925360.2
$ws.Range("I55").Value = 2942334.5
$ws.Range("J55").Value = 4408.5557
$ws.Range("K55").Value = 2942334.5
$ws.Range("L55").Value = 4408.5557
$ws.Range("M55").Value = -2942161.5
$ws.Range("N55").Value = -4754.5557
$ws.Range("H68").Value = 6529.619
$ws.Range("I68").Value = 3284.182
$ws.Range("K68").Value = 3284.182
$ws.Range("M68").Value = -2535.182
$ws.Range("H71").Value = 6529.619
$ws.Range("I71").Value = 3284.182
$ws.Range("K71").Value = 16420.91
$ws.Range("M71").Value = -12676.91
$ws.Range("H100").Value = 10750.923
$ws.Range("J100").Value = 16429
$ws.Range("L100").Value = 16429
$ws.Range("N100").Value = -17511
$ws.Range("H136").Value = 3149.6
$ws.Range("I136").Value = 2422.3333
$ws.Range("K136").Value = 7266.999899999999
$ws.Range("M136").Value = -4716.999899999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 24415.428
$ws.Range("J18").Value = 27669
$ws.Range("L18").Value = 27669
$ws.Range("N18").Value = -28015
$ws.Range("H107").Value = 1389.2222
$ws.Range("I107").Value = 1223.6154
$ws.Range("J107").Value = 1819.8
$ws.Range("K107").Value = 3670.8462
$ws.Range("L107").Value = 5459.4
$ws.Range("M107").Value = -1750.8462
$ws.Range("N107").Value = -9299.4
$ws.Range("H122").Value = 10000
$ws.Range("I122").Value = 10000
$ws.Range("K122").Value = 30000
$ws.Range("M122").Value = -27550
$ws.Range("H132").Value = 2853.577
$ws.Range("I132").Value = 2235.05
$ws.Range("K132").Value = 6705.150000000001
$ws.Range("M132").Value = -4175.150000000001
$ws.Range("H136").Value = 5407.7856
$ws.Range("I136").Value = 984.25
$ws.Range("K136").Value = 2952.75
$ws.Range("M136").Value = -402.75
